# Update cryptos list: refresh Price (D) and Volume(1h) (E) values,
# and swap two pairs of adjacent coin rows (14<->15, 16<->17) to reflect
# the new ranking order, per the "Updated cryptos list" GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "87.965.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.189.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.10%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.66%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "612.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.393"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.675"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.84%  "
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.185.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.543"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.75%  "
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.69%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.31"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.771.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "32.51"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.55%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.807.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.183.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -5.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "414.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -7.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.344.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000133"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "73.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.163"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "548.97"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.33"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.87"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.132"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.44%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "21.90"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.13%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.376"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.71%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "175.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.126"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -8.57%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.75%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.32%  "
